$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2, 3, 4 are cyclically rotated:
#   new row 2 <- old row 3
#   new row 3 <- old row 4
#   new row 4 <- old row 2 (D,J) + old row 3 (K,L,M,N,O,P,Q)
# Capture old values first (using .Value() to actually invoke the getter)
# before overwriting any of them.

$r2D = $ws.Range("D2").Value()
$r2J = $ws.Range("J2").Value()
$r2K = $ws.Range("K2").Value()
$r2L = $ws.Range("L2").Value()
$r2M = $ws.Range("M2").Value()
$r2N = $ws.Range("N2").Value()
$r2O = $ws.Range("O2").Value()
$r2P = $ws.Range("P2").Value()
$r2Q = $ws.Range("Q2").Value()

$r3D = $ws.Range("D3").Value()
$r3J = $ws.Range("J3").Value()
$r3K = $ws.Range("K3").Value()
$r3L = $ws.Range("L3").Value()
$r3M = $ws.Range("M3").Value()
$r3N = $ws.Range("N3").Value()
$r3O = $ws.Range("O3").Value()
$r3P = $ws.Range("P3").Value()
$r3Q = $ws.Range("Q3").Value()

$r4D = $ws.Range("D4").Value()
$r4J = $ws.Range("J4").Value()
$r4K = $ws.Range("K4").Value()
$r4L = $ws.Range("L4").Value()
$r4M = $ws.Range("M4").Value()
$r4N = $ws.Range("N4").Value()
$r4O = $ws.Range("O4").Value()
$r4P = $ws.Range("P4").Value()
$r4Q = $ws.Range("Q4").Value()

# Row 2 <- old row 3
$ws.Range("D2").Value = $r3D
$ws.Range("J2").Value = $r3J
$ws.Range("K2").Value = $r3K
$ws.Range("L2").Value = $r3L
$ws.Range("M2").Value = $r3M
$ws.Range("N2").Value = $r3N
$ws.Range("O2").Value = $r3O
$ws.Range("P2").Value = $r3P
$ws.Range("Q2").Value = $r3Q

# Row 3 <- old row 4
$ws.Range("D3").Value = $r4D
$ws.Range("J3").Value = $r4J
$ws.Range("K3").Value = $r4K
$ws.Range("L3").Value = $r4L
$ws.Range("M3").Value = $r4M
$ws.Range("N3").Value = $r4N
$ws.Range("O3").Value = $r4O
$ws.Range("P3").Value = $r4P
$ws.Range("Q3").Value = $r4Q

# Row 4 <- old row 2
$ws.Range("D4").Value = $r2D
$ws.Range("J4").Value = $r2J
$ws.Range("K4").Value = $r2K
$ws.Range("L4").Value = $r2L
$ws.Range("M4").Value = $r2M
$ws.Range("N4").Value = $r2N
$ws.Range("O4").Value = $r2O
$ws.Range("P4").Value = $r2P
$ws.Range("Q4").Value = $r2Q

$wb.Save()
